# Apply the United Arab Emirates_M2 data refresh:
#  - Update a handful of revised open/high/low/close values for existing rows
#  - Append three new monthly rows (258-260) with their data
#  - Dimension will be recalculated automatically by the engine on save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised values for existing rows (C,D,E,F columns share the same value) ---
$updates = @{
    226 = 1486571000000
    230 = 1496040000000
    237 = 1485900000000
    238 = 1498300000000
    239 = 1517100000000
    240 = 1563100000000
    241 = 1563400000000
    242 = 1570300000000
    243 = 1593400000000
    244 = 1567300000000
    245 = 1568100000000
    246 = 1622200000000
    247 = 1606000000000
    248 = 1627700000000
    249 = 1645600000000
    250 = 1629400000000
    251 = 1678100000000
}

foreach ($row in $updates.Keys) {
    $val = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $val
    $ws.Cells.Item($row, 4).Value = $val
    $ws.Cells.Item($row, 5).Value = $val
    $ws.Cells.Item($row, 6).Value = $val
}

# --- Append new rows 258, 259, 260 ---
$newRows = @(
    @{ Row = 258; Date = 45078.41666666666; Value = 1855306000000 },
    @{ Row = 259; Date = 45108.41666666666; Value = 1858844000000 },
    @{ Row = 260; Date = 45139.41666666666; Value = 1860330000000 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.Date
    $ws.Cells.Item($r, 2).Value = "ECONOMICS:AEM2"
    $ws.Cells.Item($r, 3).Value = $nr.Value
    $ws.Cells.Item($r, 4).Value = $nr.Value
    $ws.Cells.Item($r, 5).Value = $nr.Value
    $ws.Cells.Item($r, 6).Value = $nr.Value
    $ws.Cells.Item($r, 7).Value = 0

    # Match the date column formatting (bold, bordered, centered, date numFmt)
    # used by every other row in column A.
    $ws.Range("A257").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$ws.Range("A1").Select()
